$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47: correct the "Minor changes..." entry (B47 2 -> 1; D47 text corrected)
$ws.Range("B47").Value = 1
$ws.Range("D47").Value = "Minor changes on documentation and setup"

# Row 49 stays 2 hours, text unchanged (same string as before)
$ws.Range("B49").Value = 2
$ws.Range("D49").Value = "tc: System load estimation put to operation, validated by test case tc10"

# New rows 50 and 51 - copy A49's formatting (date number format) down first
$ws.Range("A49").Copy($ws.Range("A50"))
$ws.Range("A49").Copy($ws.Range("A51"))

$ws.Range("A50").Value = [DateTime]::FromOADate(41243)
$ws.Range("B50").Value = 3.75
$ws.Range("D50").Value = "tc05 revised, documentation of rtos.c/h extended/corrected"

$ws.Range("A51").Value = [DateTime]::FromOADate(41244)
$ws.Range("B51").Value = 1.75
$ws.Range("D51").Value = "Documentation, tc10 and doxygen"

$ws.Range("E51").Select()
